$d = $word.ActiveDocument

# Position an insertion point collapsed at the very end of the document
# (right after the "Novamente alteração" paragraph).
$endRange = $d.Range($d.Content.End, $d.Content.End)

# Insert an empty paragraph followed by a paragraph containing two runs:
# "Alteração " and "Commit & Push", using raw OOXML so the runs stay
# distinct (rather than being merged together by plain text insertion).
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r><w:t xml:space="preserve">Alteração </w:t></w:r>
            <w:r><w:t>Commit &amp; Push</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$endRange.InsertXML($xml)
